$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()), (''selector'', ''passthrough''),
                (''model'',
                 BaggingClassifier(estimator=MLPClassifier(alpha=0.01,
                                                           hidden_layer_sizes=(10,
                                                                               10,
                                                                               10),
                                                           learning_rate_init=0.0001,
                                                           max_iter=1000,
                                                           random_state=42),
                                   n_estimators=5, random_state=42))])'
$ws.Range("B2").Value = 0.6476190476190476
$ws.Range("C2").Value = '{''scaler'': MinMaxScaler(), ''model__n_estimators'': 5, ''model__estimator__solver'': ''adam'', ''model__estimator__learning_rate_init'': 0.0001, ''model__estimator__hidden_layer_sizes'': (10, 10, 10), ''model__estimator__alpha'': 0.01, ''model__estimator__activation'': ''relu''}'
$ws.Range("D2").Value = 0.5882352941176471
$ws.Range("E2").Value = '[1 0 0 1 0 0 1 1 0 1 0 0]'
$ws.Range("F2").Value = '[1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.7979285714285714
$ws.Range("I2").Value = 0.02234110707813683
$ws.Range("J2").Value = 0.5743809523809523
$ws.Range("K2").Value = 0.06506598965880907

# Row 3
$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', None), (''selector'', ''passthrough''),
                (''model'',
                 BaggingClassifier(estimator=MLPClassifier(activation=''tanh'',
                                                           alpha=0.01,
                                                           hidden_layer_sizes=(5,
                                                                               10,
                                                                               5),
                                                           learning_rate_init=1,
                                                           max_iter=1000,
                                                           random_state=42,
                                                           solver=''lbfgs''),
                                   n_estimators=50, random_state=42))])'
$ws.Range("B3").Value = 0.6476190476190476
$ws.Range("C3").Value = '{''scaler'': None, ''model__n_estimators'': 50, ''model__estimator__solver'': ''lbfgs'', ''model__estimator__learning_rate_init'': 1, ''model__estimator__hidden_layer_sizes'': (5, 10, 5), ''model__estimator__alpha'': 0.01, ''model__estimator__activation'': ''tanh''}'
$ws.Range("D3").Value = 0.5333333333333333
$ws.Range("E3").Value = '[1 0 1 0 0 0 0 1 1 0 1 1]'
$ws.Range("F3").Value = '[1 1 1 1 1 0 1 1 0 1 1 0]'
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 0.7693809523809523
$ws.Range("I3").Value = 0.01941842714124612
$ws.Range("J3").Value = 0.5555238095238095
$ws.Range("K3").Value = 0.06203243572108854

# Row 4
$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', None), (''selector'', ''passthrough''),
                (''model'',
                 BaggingClassifier(estimator=MLPClassifier(alpha=1e-05,
                                                           hidden_layer_sizes=(5,
                                                                               10,
                                                                               5),
                                                           learning_rate_init=0.0001,
                                                           max_iter=1000,
                                                           random_state=42,
                                                           solver=''sgd''),
                                   n_estimators=50, random_state=42))])'
$ws.Range("B4").Value = 0.6
$ws.Range("C4").Value = '{''scaler'': None, ''model__n_estimators'': 50, ''model__estimator__solver'': ''sgd'', ''model__estimator__learning_rate_init'': 0.0001, ''model__estimator__hidden_layer_sizes'': (5, 10, 5), ''model__estimator__alpha'': 1e-05, ''model__estimator__activation'': ''relu''}'
$ws.Range("D4").Value = 0.8
$ws.Range("E4").Value = '[1 0 1 1 1 1 0 1 0 1 0 1]'
$ws.Range("F4").Value = '[1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.742761904761905
$ws.Range("I4").Value = 0.02184327197518508
$ws.Range("J4").Value = 0.5434285714285714
$ws.Range("K4").Value = 0.06465148221885886
